# Update "想去人数" (want-to-go count) figures in column F for the
# 展览 and 全部类型 sheets, per the latest data pull.

$wb = $excel.ActiveWorkbook

# Row -> [old value, new value] for column F on the affected sheets.
$updates = @{
    3  = 298
    5  = 606
    7  = 2092
    10 = 4573
    12 = 290
    15 = 143
    18 = 90
    19 = 3473
    21 = 554
    25 = 98
    29 = 211
    31 = 694
    32 = 2107
    33 = 397
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
